$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.27

$ws.Range("B3").Value = 1.48
$ws.Range("E3").Value = 1.29
$ws.Range("F3").Value = 1.22

$ws.Range("B4").Value = 1.47
$ws.Range("C4").Value = 1.45
$ws.Range("D4").Value = 1.32
$ws.Range("E4").Value = 1.23

$ws.Range("C5").Value = 1.35
$ws.Range("D5").Value = 1.35
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 0.78

$ws.Range("D6").Value = 1.54
$ws.Range("E6").Value = 1.32

$ws.Range("E7").Value = 1.87
$ws.Range("G7").Value = 1.18
